$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fecha" (timestamp) values pushed down the log on this update run,
# per "Actualizar 02-05-2021 13-20-53":
#   rows 2-15  (newest block)      -> new check timestamp
#   rows 16-29 (was rows 2-15)     -> previous newest timestamp, shifted down
#   rows 30-37 (was rows 16-23)    -> previous timestamp, shifted down
$newest = 44232.55611753627
$shift1 = 44232.53506424768
$shift2 = 44232.51400987268

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift1
}

for ($r = 30; $r -le 37; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift2
}
